$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing description text in G28 (Unet -> PSPNet on 2nd line) ---
$ws.Range("G28").Value = "1. PSPNet_starter nb`n2.  checked other PSPNet implementation differences"

# --- Add new row 29 (Sno 28) entry ---
# Copy formatting from row 28 (same style pattern) into row 29 first
$ws.Range("A28:G28").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 44738
$ws.Range("C29").Value = 0.52083333333333337
$ws.Range("D29").Value = 0.625
$ws.Range("E29").Formula = "=D29-C29"
$ws.Range("F29").Value = "Code"
$ws.Range("G29").Value = "1. PSPNet_starter nb first phase completed`n2. PSPNet_resnet50_baseline model run without aux_loss for 10 epochs and save model output, video`n3. PSPNet_resnet50_aux model run with aux_loss for 10 epochs and save model output, video"
$ws.Rows(29).RowHeight = 75

# --- Update sheet view: scroll down and move selection to G30 ---
$ws.Range("G30").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1

# --- Recalculate so the Total Hours sum (E33) reflects the new row ---
$excel.CalculateFull() | Out-Null
